# Auto-generated Excel COM-interop edit script
# Applies the diff: regenerate results/multi_city/AP/desc_stats.xlsx by ignoring CIs for health (Rel #120)

$wb = $excel.ActiveWorkbook
$wsPM = $wb.Worksheets.Item("PM_2.5_emission_inventory")
$wsSummary = $wb.Worksheets.Item("summary_stats_PM2.5_CO2")

# --- Sheet "PM_2.5_emission_inventory": update cape_town column (X) rows 2-7 ---
# These cells are stored as text (e.g. "39.6"), not numbers, so force Text type
# using an apostrophe prefix to prevent Excel auto-converting them to numeric values.
$wsPM.Range("X2").Value = "'39.6"
$wsPM.Range("X3").Value = "'21.6"
$wsPM.Range("X4").Value = "'20.5"
$wsPM.Range("X5").Value = "'0.7"
$wsPM.Range("X6").Value = "'0"
$wsPM.Range("X7").Value = "'17.6"

# --- Sheet "summary_stats_PM2.5_CO2": update descriptive-stat rows (numeric columns B-G) ---
# row 107
$wsSummary.Range("B107").Value = 33.3
$wsSummary.Range("C107").Value = 32.4
$wsSummary.Range("D107").Value = 32.4
$wsSummary.Range("E107").Value = 34.9
$wsSummary.Range("F107").Value = 32.38

# row 108
$wsSummary.Range("B108").Value = 33.1
$wsSummary.Range("C108").Value = 32
$wsSummary.Range("D108").Value = 32
$wsSummary.Range("E108").Value = 35.2
$wsSummary.Range("F108").Value = 32.04
$wsSummary.Range("G108").Value = -0.34

# row 109
$wsSummary.Range("B109").Value = 36.6
$wsSummary.Range("C109").Value = 35.6
$wsSummary.Range("D109").Value = 35.6
$wsSummary.Range("E109").Value = 38.4
$wsSummary.Range("F109").Value = 35.63
$wsSummary.Range("G109").Value = 3.25

# row 110
$wsSummary.Range("B110").Value = 32.9
$wsSummary.Range("C110").Value = 32
$wsSummary.Range("D110").Value = 32
$wsSummary.Range("E110").Value = 34.5
$wsSummary.Range("F110").Value = 32.03
$wsSummary.Range("G110").Value = -0.35

# row 111
$wsSummary.Range("B111").Value = 32.4
$wsSummary.Range("C111").Value = 31.6
$wsSummary.Range("D111").Value = 31.6
$wsSummary.Range("E111").Value = 34
$wsSummary.Range("F111").Value = 31.56
$wsSummary.Range("G111").Value = -0.82

# row 112
$wsSummary.Range("B112").Value = 15.9
$wsSummary.Range("C112").Value = 14.9
$wsSummary.Range("D112").Value = 15.4
$wsSummary.Range("E112").Value = 17
$wsSummary.Range("F112").Value = 14.86

# row 113
$wsSummary.Range("B113").Value = 15.8
$wsSummary.Range("C113").Value = 14.7
$wsSummary.Range("D113").Value = 15.4
$wsSummary.Range("E113").Value = 17.1
$wsSummary.Range("F113").Value = 14.72
$wsSummary.Range("G113").Value = -0.14

# row 114
$wsSummary.Range("B114").Value = 15.8
$wsSummary.Range("C114").Value = 14.8
$wsSummary.Range("D114").Value = 15.3
$wsSummary.Range("E114").Value = 17
$wsSummary.Range("F114").Value = 14.78
$wsSummary.Range("G114").Value = -0.08

# row 115
$wsSummary.Range("B115").Value = 15.8
$wsSummary.Range("C115").Value = 14.8
$wsSummary.Range("D115").Value = 15.4
$wsSummary.Range("E115").Value = 16.9
$wsSummary.Range("F115").Value = 14.83
$wsSummary.Range("G115").Value = -0.03

# row 116
$wsSummary.Range("B116").Value = 16.1
$wsSummary.Range("C116").Value = 15.1
$wsSummary.Range("D116").Value = 15.6
$wsSummary.Range("E116").Value = 17.3
$wsSummary.Range("F116").Value = 15.06
$wsSummary.Range("G116").Value = 0.2

# row 117
$wsSummary.Range("B117").Value = 109.7
$wsSummary.Range("C117").Value = 105
$wsSummary.Range("D117").Value = 107.4
$wsSummary.Range("E117").Value = 118
$wsSummary.Range("F117").Value = 105

# row 118
$wsSummary.Range("B118").Value = 108.6
$wsSummary.Range("C118").Value = 103.5
$wsSummary.Range("D118").Value = 105.8
$wsSummary.Range("E118").Value = 117.9
$wsSummary.Range("F118").Value = 103.49
$wsSummary.Range("G118").Value = -1.51

# row 119
$wsSummary.Range("B119").Value = 111.4
$wsSummary.Range("C119").Value = 106.7
$wsSummary.Range("D119").Value = 109.1
$wsSummary.Range("E119").Value = 120.1
$wsSummary.Range("F119").Value = 106.69
$wsSummary.Range("G119").Value = 1.69

# row 120
$wsSummary.Range("B120").Value = 108.4
$wsSummary.Range("C120").Value = 103.8
$wsSummary.Range("D120").Value = 106.1
$wsSummary.Range("E120").Value = 116.7
$wsSummary.Range("F120").Value = 103.78
$wsSummary.Range("G120").Value = -1.22

# row 121
$wsSummary.Range("B121").Value = 108.2
$wsSummary.Range("C121").Value = 103.7
$wsSummary.Range("D121").Value = 105.9
$wsSummary.Range("E121").Value = 116.4
$wsSummary.Range("F121").Value = 103.67
$wsSummary.Range("G121").Value = -1.33

# row 122
$wsSummary.Range("B122").Value = 34.1
$wsSummary.Range("C122").Value = 32.8
$wsSummary.Range("D122").Value = 33.6
$wsSummary.Range("E122").Value = 36.2
$wsSummary.Range("F122").Value = 32.75

# row 123
$wsSummary.Range("B123").Value = 34.1
$wsSummary.Range("C123").Value = 32.5
$wsSummary.Range("D123").Value = 33.3
$wsSummary.Range("E123").Value = 36.5
$wsSummary.Range("F123").Value = 32.5
$wsSummary.Range("G123").Value = -0.25

# row 124
$wsSummary.Range("B124").Value = 39.4
$wsSummary.Range("C124").Value = 37.7
$wsSummary.Range("D124").Value = 38.7
$wsSummary.Range("E124").Value = 41.7
$wsSummary.Range("F124").Value = 37.71
$wsSummary.Range("G124").Value = 4.96

# row 125
$wsSummary.Range("B125").Value = 33.6
$wsSummary.Range("C125").Value = 32.2
$wsSummary.Range("D125").Value = 33.1
$wsSummary.Range("E125").Value = 35.6
$wsSummary.Range("F125").Value = 32.23
$wsSummary.Range("G125").Value = -0.52

# row 126
$wsSummary.Range("B126").Value = 34
$wsSummary.Range("C126").Value = 32.6
$wsSummary.Range("D126").Value = 33.4
$wsSummary.Range("E126").Value = 35.9
$wsSummary.Range("F126").Value = 32.58
$wsSummary.Range("G126").Value = -0.17

# row 137
$wsSummary.Range("B137").Value = 14.7
$wsSummary.Range("C137").Value = 14
$wsSummary.Range("D137").Value = 14.5
$wsSummary.Range("E137").Value = 15.5
$wsSummary.Range("F137").Value = 14

# row 138
$wsSummary.Range("B138").Value = 14.7
$wsSummary.Range("C138").Value = 13.9
$wsSummary.Range("D138").Value = 14.5
$wsSummary.Range("E138").Value = 15.7
$wsSummary.Range("F138").Value = 13.93
$wsSummary.Range("G138").Value = -0.07000000000000001

# row 139
$wsSummary.Range("B139").Value = 15.3
$wsSummary.Range("C139").Value = 14.6
$wsSummary.Range("D139").Value = 15.1
$wsSummary.Range("E139").Value = 16.1
$wsSummary.Range("F139").Value = 14.6
$wsSummary.Range("G139").Value = 0.6

# row 140
$wsSummary.Range("B140").Value = 14.2
$wsSummary.Range("C140").Value = 13.6
$wsSummary.Range("D140").Value = 14.1
$wsSummary.Range("E140").Value = 15
$wsSummary.Range("F140").Value = 13.62
$wsSummary.Range("G140").Value = -0.38

# row 141
$wsSummary.Range("B141").Value = 14.3
$wsSummary.Range("C141").Value = 13.7
$wsSummary.Range("D141").Value = 14.2
$wsSummary.Range("E141").Value = 15.1
$wsSummary.Range("F141").Value = 13.7
$wsSummary.Range("G141").Value = -0.3

